$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue 2 4 "27.467.37"
Set-TextValue 2 5 "  -0.28%  "
Set-TextValue 3 4 "1.615.73"
Set-TextValue 3 5 "  -1.50%  "
Set-TextValue 4 5 "  +0.17%  "
Set-TextValue 5 5 "  -0.87%  "
Set-TextValue 6 5 "  -2.00%  "
Set-TextValue 7 5 "  +0.18%  "
Set-TextValue 8 4 "22.91"
Set-TextValue 8 5 "  -0.24%  "
Set-TextValue 9 5 "  +1.01%  "
Set-TextValue 10 5 "  -0.03%  "
Set-TextValue 11 5 "  -0.55%  "
Set-TextValue 12 4 "1.843.56"
Set-TextValue 12 5 "  -1.55%  "
Set-TextValue 13 4 "1.618.41"
Set-TextValue 13 5 "  -1.49%  "
Set-TextValue 14 5 "  -0.12%  "
Set-TextValue 15 4 "0.550"
Set-TextValue 15 5 "  -2.61%  "
Set-TextValue 16 4 "64.29"
Set-TextValue 16 5 "  +0.12%  "
Set-TextValue 17 4 "27.472.75"
Set-TextValue 17 5 "  -0.31%  "
Set-TextValue 18 4 "227.84"
Set-TextValue 18 5 "  -0.62%  "
Set-TextValue 20 4 "7.54"
Set-TextValue 20 5 "  -2.32%  "
Set-TextValue 21 5 "  +0.13%  "
Set-TextValue 22 4 "4.28"
Set-TextValue 22 5 "  -0.76%  "
Set-TextValue 23 4 "9.90"
Set-TextValue 23 5 "  -0.15%  "
Set-TextValue 24 4 "2.08"
Set-TextValue 24 5 "  +6.52%  "
Set-TextValue 25 4 "149.16"
Set-TextValue 25 5 "  -0.29%  "
Set-TextValue 26 5 "  -1.20%  "
Set-TextValue 27 2 "Cosmos"
Set-TextValue 27 3 "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue 27 4 "6.82"
Set-TextValue 27 5 "  -2.11%  "
Set-TextValue 28 2 "BinanceUSD"
Set-TextValue 28 3 "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue 28 4 "1.00"
Set-TextValue 28 5 "  +0.15%  "
Set-TextValue 29 4 "15.56"
Set-TextValue 29 5 "  -0.27%  "
Set-TextValue 30 5 "  -0.88%  "
Set-TextValue 31 5 "  -1.29%  "
Set-TextValue 32 5 "  -0.47%  "
Set-TextValue 33 4 "1.444.38"
Set-TextValue 33 5 "  +1.20%  "
Set-TextValue 34 4 "3.04"
Set-TextValue 34 5 "  -3.78%  "
Set-TextValue 35 4 "1.53"
Set-TextValue 35 5 "  -3.28%  "
Set-TextValue 36 5 "  -0.17%  "
Set-TextValue 37 4 "0.940"
Set-TextValue 37 5 "  +5.11%  "
Set-TextValue 38 5 "  -1.62%  "
Set-TextValue 39 5 "  +0.09%  "
Set-TextValue 40 4 "0.861"
Set-TextValue 40 5 "  -2.11%  "
Set-TextValue 41 4 "69.19"
Set-TextValue 41 5 "  +6.52%  "
Set-TextValue 42 5 "  +0.14%  "
Set-TextValue 43 5 "  -1.94%  "
Set-TextValue 44 5 "  +0.22%  "
Set-TextValue 45 4 "5.38"
Set-TextValue 45 5 "  -2.53%  "
Set-TextValue 46 4 "2.21"
Set-TextValue 46 5 "  -2.13%  "
Set-TextValue 47 4 "1.756.31"
Set-TextValue 47 5 "  -1.42%  "
Set-TextValue 48 5 "  -0.06%  "
Set-TextValue 49 4 "86.11"
Set-TextValue 49 5 "  -0.02%  "
Set-TextValue 50 5 "  -1.34%  "
Set-TextValue 51 4 "0.0987"
Set-TextValue 51 5 "  -0.10%  "
